$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new date and a longer note (now 3 lines -> taller row)
$ws.Range("D6").Value = "Poznámka na `noveľa`nviac riadkov"
$ws.Range("C6").Value = [datetime]"2019-09-04"
$ws.Rows.Item(6).RowHeight = 75

# Row 5: new contact block (name / email) and a new date
$ws.Range("B5").Value = "jean@gmail.com`njakub.roncak@gmail.kroň`nelf@rene.sk"
$ws.Range("A5").Value = "Erik Laďnakghj`nJéáň Réňo`nElf Reňe"
$ws.Range("C5").Value = [datetime]"2020-08-04"

# Active cell / selection moves to C5
$ws.Range("C5").Select() | Out-Null

Write-Output "done"
